# Update latest output (run 89)

$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" ---
$ws1 = $wb.Worksheets.Item("Schedule")

$ws1.Range("A2").Value = 46040.27083333334
$ws1.Range("B2").Value = 46040.89583333334
$ws1.Range("C2").Value = 15
$ws1.Range("D2").Value = 56.7
$ws1.Range("E2").Value = 196.736553
$ws1.Range("F2").Value = 3.469780476190477

$ws1.Range("A3").Value = 46040.95833333334
$ws1.Range("B3").Value = 46041.14583333334
$ws1.Range("E3").Value = 328.78685775
$ws1.Range("F3").Value = 19.32903337742505

$ws1.Range("A4").Value = 46041.3125
$ws1.Range("C4").Value = 8.5
$ws1.Range("D4").Value = 32.13
$ws1.Range("E4").Value = 31.82877749999998
$ws1.Range("F4").Value = 0.9906248832866477

# --- Sheet "Detailed" ---
$ws2 = $wb.Worksheets.Item("Detailed")

$ws2.Range("E14").Value = "OFF"

$ws2.Range("B43").Value = 55.07604
$ws2.Range("E43").Value = "ON"

$ws2.Range("B44").Value = 56.42913
$ws2.Range("E44").Value = "ON"

$ws2.Range("C45").Value = "historical"
$ws2.Range("C46").Value = "historical"

$ws2.Range("E47").Value = "OFF"

$ws2.Range("B48").Value = 47.30732
$ws2.Range("B49").Value = 47.21051
$ws2.Range("B50").Value = 47.34256
$ws2.Range("B51").Value = 36.2
$ws2.Range("B52").Value = 36.2
$ws2.Range("B53").Value = 28.97997
$ws2.Range("B54").Value = 28.81854
$ws2.Range("B55").Value = 29.27843
$ws2.Range("B56").Value = 35.87996
$ws2.Range("E56").Value = "ON"

$ws2.Range("B57").Value = 36.2

$ws2.Range("B59").Value = 58.70125
$ws2.Range("B60").Value = 58.59489
$ws2.Range("B61").Value = 60.20729

$ws2.Range("E64").Value = "OFF"

$ws2.Range("B65").Value = 35.71604
$ws2.Range("B66").Value = 28.53116

$ws2.Range("B68").Value = 0.62379
$ws2.Range("B69").Value = -4.92073
$ws2.Range("B70").Value = -5.01
$ws2.Range("B71").Value = -5.01
$ws2.Range("B72").Value = -5.74313
$ws2.Range("B73").Value = -2.83936
$ws2.Range("B74").Value = -5.45907
$ws2.Range("B75").Value = -5.75885
$ws2.Range("B76").Value = -6.10743
$ws2.Range("B77").Value = -5.99642
$ws2.Range("B78").Value = -4.89752
$ws2.Range("B79").Value = -2.55361
$ws2.Range("B80").Value = 0.00001
$ws2.Range("B81").Value = 0.00002
$ws2.Range("B82").Value = 0.00025
$ws2.Range("B83").Value = -5.25437
$ws2.Range("B84").Value = -2.52116
$ws2.Range("B85").Value = 0.854
$ws2.Range("B86").Value = 12.22181
$ws2.Range("B87").Value = 53.21494
$ws2.Range("B88").Value = 60.45778
$ws2.Range("B89").Value = 65
$ws2.Range("B90").Value = 68.21548
$ws2.Range("B91").Value = 60.30289
$ws2.Range("B92").Value = 65

$ws2.Range("B94").Value = 62.11894
$ws2.Range("B95").Value = 58.97459
$ws2.Range("B96").Value = 58.14024
$ws2.Range("B97").Value = 62.03932
